$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E ("Group"), shifting E:V right to
# F:W. Shifting right copies column D's formatting into both the new
# column E and the pushed-over column F (matches Excel's default
# insert-column behaviour).
$ws.Columns("E").Insert(-4161)

# The old "Group" header (now sitting in F3) is renamed to "Program Group".
$ws.Range("F3").Value = "Program Group"

# The newly inserted column gets the new "Collection" header.
$ws.Range("E3").Value = "Collection"

# Widen the "Collection" and "Program Group" columns to fit their longer
# text; column D keeps its original width. (25.67 is the COM ColumnWidth
# input that round-trips to the target stored width of 26.5703125.)
$ws.Range("E1").ColumnWidth = 25.67
$ws.Range("F1").ColumnWidth = 25.67
